$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.306379666666667
$ws.Range("H2").Value = 3.919139
$ws.Range("I2").Value = 0.4034923136874173
$ws.Range("J2").Value = 0.4034923136874172
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.39487933333333
$ws.Range("N2").Value = 43.184638
$ws.Range("O2").Value = 0.6893176795959761
$ws.Range("P2").Value = 0.6893176795959762
$ws.Range("Q2").Value = 18.80517766518689
$ws.Range("R2").Value = 169.246598986682
$ws.Range("S2").Value = 0.2781343854058222
$ws.Range("T2").Value = 0.2781343854058222
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.306379666666667
$ws.Range("H3").Value = 3.919139
$ws.Range("I3").Value = 0.4034923136874173
$ws.Range("J3").Value = 0.4034923136874172
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.526052333333333
$ws.Range("N3").Value = 16.578157
$ws.Range("O3").Value = 0.264622264871545
$ws.Range("P3").Value = 0.2646222648715451
$ws.Range("Q3").Value = 7.219122405202556
$ws.Range("R3").Value = 64.97210164682301
$ws.Range("S3").Value = 0.1067730499062243
$ws.Range("T3").Value = 0.1067730499062243
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.306379666666667
$ws.Range("H4").Value = 3.919139
$ws.Range("I4").Value = 0.4034923136874173
$ws.Range("J4").Value = 0.4034923136874172
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.9618626666666668
$ws.Range("N4").Value = 2.885588
$ws.Range("O4").Value = 0.04606005553247879
$ws.Range("P4").Value = 0.04606005553247879
$ws.Range("Q4").Value = 1.256557829859111
$ws.Range("R4").Value = 11.309020468732
$ws.Range("S4").Value = 0.01858487837537079
$ws.Range("T4").Value = 0.01858487837537079
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.072562
$ws.Range("H5").Value = 3.217686
$ws.Range("I5").Value = 0.3312746929515923
$ws.Range("J5").Value = 0.3312746929515923
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.39487933333333
$ws.Range("N5").Value = 43.184638
$ws.Range("O5").Value = 0.6893176795959761
$ws.Range("P5").Value = 0.6893176795959762
$ws.Range("Q5").Value = 15.43940056751866
$ws.Range("R5").Value = 138.954605107668
$ws.Range("S5").Value = 0.2283535026542611
$ws.Range("T5").Value = 0.2283535026542611
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.072562
$ws.Range("H6").Value = 3.217686
$ws.Range("I6").Value = 0.3312746929515923
$ws.Range("J6").Value = 0.3312746929515923
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.526052333333333
$ws.Range("N6").Value = 16.578157
$ws.Range("O6").Value = 0.264622264871545
$ws.Range("P6").Value = 0.2646222648715451
$ws.Range("Q6").Value = 5.927033742744666
$ws.Range("R6").Value = 53.34330368470199
$ws.Range("S6").Value = 0.087662659543476
$ws.Range("T6").Value = 0.08766265954347605
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.072562
$ws.Range("H7").Value = 3.217686
$ws.Range("I7").Value = 0.3312746929515923
$ws.Range("J7").Value = 0.3312746929515923
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9618626666666668
$ws.Range("N7").Value = 2.885588
$ws.Range("O7").Value = 0.04606005553247879
$ws.Range("P7").Value = 0.04606005553247879
$ws.Range("Q7").Value = 1.031657345485333
$ws.Range("R7").Value = 9.284916109368
$ws.Range("S7").Value = 0.0152585307538552
$ws.Range("T7").Value = 0.01525853075385521
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.8587400000000001
$ws.Range("H8").Value = 2.57622
$ws.Range("I8").Value = 0.2652329933609903
$ws.Range("J8").Value = 0.2652329933609903
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.39487933333333
$ws.Range("N8").Value = 43.184638
$ws.Range("O8").Value = 0.6893176795959761
$ws.Range("P8").Value = 0.6893176795959762
$ws.Range("Q8").Value = 12.36145867870667
$ws.Range("R8").Value = 111.25312810836
$ws.Range("S8").Value = 0.1828297915358928
$ws.Range("T8").Value = 0.1828297915358928
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.8587400000000001
$ws.Range("H9").Value = 2.57622
$ws.Range("I9").Value = 0.2652329933609903
$ws.Range("J9").Value = 0.2652329933609903
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.526052333333333
$ws.Range("N9").Value = 16.578157
$ws.Range("O9").Value = 0.264622264871545
$ws.Range("P9").Value = 0.2646222648715451
$ws.Range("Q9").Value = 4.745442180726667
$ws.Range("R9").Value = 42.70897962654001
$ws.Range("S9").Value = 0.07018655542184472
$ws.Range("T9").Value = 0.07018655542184474
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8587400000000001
$ws.Range("H10").Value = 2.57622
$ws.Range("I10").Value = 0.2652329933609903
$ws.Range("J10").Value = 0.2652329933609903
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.9618626666666668
$ws.Range("N10").Value = 2.885588
$ws.Range("O10").Value = 0.04606005553247879
$ws.Range("P10").Value = 0.04606005553247879
$ws.Range("Q10").Value = 0.8259899463733335
$ws.Range("R10").Value = 7.433909517360001
$ws.Range("S10").Value = 0.01221664640325279
$ws.Range("T10").Value = 0.01221664640325279
